# Updates the cryptos list (Price / Volume(1h) columns, plus a Bittensor/Maker
# row swap) to reflect the refreshed figures from the GitHub Actions scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force-write a value as TEXT so Excel does not reinterpret
# numeric-looking strings (e.g. '6.73') as real numbers, while keeping
# the cell's original (default) style untouched.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '67.857.78'
$ws.Range("E2").Value = '  +4.78%  '
$ws.Range("D3").Value = '3.266.18'
$ws.Range("E3").Value = '  +4.79%  '
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue $ws.Range("D5") '580.38'
$ws.Range("E5").Value = '  +2.83%  '
Set-TextValue $ws.Range("D6") '182.11'
$ws.Range("E6").Value = '  +8.84%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("D9").Value = '3.265.36'
$ws.Range("E9").Value = '  +4.80%  '
$ws.Range("E10").Value = '  +9.90%  '
Set-TextValue $ws.Range("D11") '6.73'
$ws.Range("E11").Value = '  +3.50%  '
Set-TextValue $ws.Range("D12") '0.416'
$ws.Range("E12").Value = '  +8.19%  '
$ws.Range("D13").Value = '3.832.80'
$ws.Range("E13").Value = '  +4.89%  '
$ws.Range("E14").Value = '  +1.58%  '
Set-TextValue $ws.Range("D15") '28.44'
$ws.Range("E15").Value = '  +6.75%  '
$ws.Range("D16").Value = '67.808.13'
$ws.Range("E16").Value = '  +4.73%  '
$ws.Range("E17").Value = '  +5.33%  '
$ws.Range("D18").Value = '3.264.52'
$ws.Range("E18").Value = '  +4.86%  '
$ws.Range("E19").Value = '  +4.51%  '
Set-TextValue $ws.Range("D20") '13.53'
$ws.Range("E20").Value = '  +7.76%  '
Set-TextValue $ws.Range("D21") '375.50'
$ws.Range("E22").Value = '  +7.63%  '
$ws.Range("E23").Value = '  +0.22%  '
Set-TextValue $ws.Range("D24") '71.02'
$ws.Range("E24").Value = '  +4.28%  '
$ws.Range("E25").Value = '  +5.12%  '
$ws.Range("E26").Value = '  +8.99%  '
$ws.Range("E27").Value = '  +1.63%  '
Set-TextValue $ws.Range("D28") '0.182'
$ws.Range("E28").Value = '  +4.05%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  +4.89%  '
Set-TextValue $ws.Range("D31") '5.68'
$ws.Range("E31").Value = '  +9.66%  '
Set-TextValue $ws.Range("D32") '22.74'
$ws.Range("E32").Value = '  +5.86%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +9.48%  '
Set-TextValue $ws.Range("D35") '6.92'
$ws.Range("E35").Value = '  +6.77%  '
Set-TextValue $ws.Range("D36") '163.73'
$ws.Range("E36").Value = '  +3.34%  '
$ws.Range("E37").Value = '  +7.25%  '
Set-TextValue $ws.Range("D38") '0.851'
$ws.Range("E38").Value = '  +4.41%  '
$ws.Range("E39").Value = '  +7.48%  '
Set-TextValue $ws.Range("D40") '6.84'
$ws.Range("E40").Value = '  +14.14%  '
Set-TextValue $ws.Range("D41") '26.76'
$ws.Range("E41").Value = '  +3.42%  '
Set-TextValue $ws.Range("D42") '4.65'
$ws.Range("E42").Value = '  +13.85%  '
Set-TextValue $ws.Range("D43") '2.60'
$ws.Range("E43").Value = '  +9.25%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.701.23'
$ws.Range("E44").Value = '  +3.51%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D45") '352.54'
$ws.Range("E45").Value = '  +11.52%  '
Set-TextValue $ws.Range("D46") '25.39'
$ws.Range("E46").Value = '  +8.98%  '
Set-TextValue $ws.Range("D47") '40.89'
$ws.Range("E47").Value = '  +4.16%  '
Set-TextValue $ws.Range("D48") '0.0680'
$ws.Range("E48").Value = '  +6.19%  '
$ws.Range("E49").Value = '  +5.09%  '
$ws.Range("E50").Value = '  +8.45%  '
$ws.Range("E51").Value = '  +1.75%  '
